# Macroferia Regional de Talca - Arandano (blue): weekly update.
# A new week's price row is inserted at row 22 (pushing the existing
# rows 22:70 down to 23:71), and the sheet's used-range dimension grows
# from A1:T70 to A1:T71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22; everything below shifts down one row.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with this week's record.
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "Macroferia Regional de Talca"
$ws.Range("C22").Value = "Maule"
$ws.Range("D22").Value = 44659
$ws.Range("E22").Value = 7
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100101
$ws.Range("H22").Value = "Berries"
$ws.Range("I22").Value = 100101001
$ws.Range("J22").Value = "Arándano (blue)"
$ws.Range("K22").Value = "Sin especificar"
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 3000
$ws.Range("O22").Value = 3000
$ws.Range("P22").Value = 3000
$ws.Range("Q22").Value = "$/bandeja 2 kilos"
$ws.Range("R22").Value = "Provincia de Linares"
$ws.Range("S22").Value = 1500
$ws.Range("T22").Value = 2
